$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D:D").Insert()

# Copy the number formatting from column E (the old column D, now shifted) into the
# new column D so the new column inherits the same date/number styles without
# introducing new style entries.
$ws.Range("E7:E102").Copy() | Out-Null
$ws.Range("D7:D102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate the refreshed financial data (new latest-year column D plus restated
# historical columns E:K) row by row.
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43100
$ws.Range("F7").Value2 = 42735
$ws.Range("G7").Value2 = 42369
$ws.Range("H7").Value2 = 42004
$ws.Range("I7").Value2 = 41639
$ws.Range("J7").Value2 = 41274
$ws.Range("K7").Value2 = 40908
$ws.Range("D8").Value2 = 84885300
$ws.Range("E8").Value2 = 84089800
$ws.Range("F8").Value2 = 82011900
$ws.Range("G8").Value2 = 77673100
$ws.Range("H8").Value2 = 70301600
$ws.Range("I8").Value2 = 67467500
$ws.Range("J8").Value2 = 65265000
$ws.Range("K8").Value2 = 68847500
$ws.Range("D9").Value2 = 40101000
$ws.Range("E9").Value2 = 40221100
$ws.Range("F9").Value2 = 39251700
$ws.Range("G9").Value2 = 47095500
$ws.Range("H9").Value2 = 43240400
$ws.Range("I9").Value2 = 40677700
$ws.Range("J9").Value2 = 38434900
$ws.Range("K9").Value2 = 39774600
$ws.Range("D10").Value2 = 44784200
$ws.Range("E10").Value2 = 43868700
$ws.Range("F10").Value2 = 42760200
$ws.Range("G10").Value2 = 30577600
$ws.Range("H10").Value2 = 27061300
$ws.Range("I10").Value2 = 26789800
$ws.Range("J10").Value2 = 26830100
$ws.Range("K10").Value2 = 29072900
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("F12").Value2 = "NA"
$ws.Range("G12").Value2 = "NA"
$ws.Range("H12").Value2 = "NA"
$ws.Range("I12").Value2 = "NA"
$ws.Range("J12").Value2 = "NA"
$ws.Range("K12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("F13").Value2 = 0
$ws.Range("G13").Value2 = 0
$ws.Range("H13").Value2 = 0
$ws.Range("I13").Value2 = 0
$ws.Range("J13").Value2 = 0
$ws.Range("K13").Value2 = 0
$ws.Range("D14").Value2 = 1350900
$ws.Range("E14").Value2 = 550900
$ws.Range("F14").Value2 = 757300
$ws.Range("G14").Value2 = 166100
$ws.Range("H14").Value2 = 133500
$ws.Range("I14").Value2 = 917800
$ws.Range("J14").Value2 = 24730900
$ws.Range("K14").Value2 = 3917000
$ws.Range("D15").Value2 = 14712700
$ws.Range("E15").Value2 = 13946300
$ws.Range("F15").Value2 = 14235800
$ws.Range("G15").Value2 = "NA"
$ws.Range("H15").Value2 = "NA"
$ws.Range("I15").Value2 = "NA"
$ws.Range("J15").Value2 = "NA"
$ws.Range("K15").Value2 = "NA"
$ws.Range("D17").Value2 = 75880200
$ws.Range("E17").Value2 = 73562200
$ws.Range("F17").Value2 = 71729900
$ws.Range("G17").Value2 = 69787800
$ws.Range("H17").Value2 = 62170600
$ws.Range("I17").Value2 = 61936100
$ws.Range("J17").Value2 = 69710400
$ws.Range("K17").Value2 = 62290600
$ws.Range("D18").Value2 = 9005100
$ws.Range("E18").Value2 = 10527600
$ws.Range("F18").Value2 = 10281900
$ws.Range("G18").Value2 = 7885300
$ws.Range("H18").Value2 = 8131100
$ws.Range("I18").Value2 = 5531400
$ws.Range("J18").Value2 = -4445300
$ws.Range("K18").Value2 = 6556900
$ws.Range("D20").Value2 = -874000
$ws.Range("E20").Value2 = -2100400
$ws.Range("F20").Value2 = -2134000
$ws.Range("G20").Value2 = 402800
$ws.Range("H20").Value2 = -260300
$ws.Range("I20").Value2 = -462300
$ws.Range("J20").Value2 = -81900
$ws.Range("K20").Value2 = 30500
$ws.Range("D21").Value2 = 23692300
$ws.Range("E21").Value2 = 24832000
$ws.Range("F21").Value2 = 23196200
$ws.Range("G21").Value2 = 21064600
$ws.Range("H21").Value2 = 19763200
$ws.Range("I21").Value2 = 17332800
$ws.Range("J21").Value2 = 20167600
$ws.Range("K21").Value2 = 23544100
$ws.Range("D22").Value2 = 2349400
$ws.Range("E22").Value2 = 2824000
$ws.Range("F22").Value2 = 3046200
$ws.Range("G22").Value2 = 2927300
$ws.Range("H22").Value2 = 2990100
$ws.Range("I22").Value2 = 2681600
$ws.Range("J22").Value2 = 2624300
$ws.Range("K22").Value2 = 3043700
$ws.Range("D23").Value2 = 5781600
$ws.Range("E23").Value2 = 5603200
$ws.Range("F23").Value2 = 5101700
$ws.Range("G23").Value2 = 5360900
$ws.Range("H23").Value2 = 4880700
$ws.Range("I23").Value2 = 2387600
$ws.Range("J23").Value2 = -7151600
$ws.Range("K23").Value2 = 3543700
$ws.Range("D24").Value2 = 2046500
$ws.Range("E24").Value2 = -626100
$ws.Range("F24").Value2 = 1619000
$ws.Range("G24").Value2 = 1431700
$ws.Range("H24").Value2 = 1240900
$ws.Range("I24").Value2 = 1036700
$ws.Range("J24").Value2 = -1700900
$ws.Range("K24").Value2 = 2757300
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("F25").Value2 = 0
$ws.Range("G25").Value2 = 0
$ws.Range("H25").Value2 = 0
$ws.Range("I25").Value2 = 0
$ws.Range("J25").Value2 = 0
$ws.Range("K25").Value2 = 0
$ws.Range("D26").Value2 = 3735100
$ws.Range("E26").Value2 = 6229300
$ws.Range("F26").Value2 = 3482700
$ws.Range("G26").Value2 = 3929200
$ws.Range("H26").Value2 = 3639700
$ws.Range("I26").Value2 = 1350900
$ws.Range("J26").Value2 = -5450600
$ws.Range("K26").Value2 = 786500
$ws.Range("D27").Value2 = 2430200
$ws.Range("E27").Value2 = 3883200
$ws.Range("F27").Value2 = 3001300
$ws.Range("G27").Value2 = 3651000
$ws.Range("H27").Value2 = 3280700
$ws.Range("I27").Value2 = 1043500
$ws.Range("J27").Value2 = -6006000
$ws.Range("K27").Value2 = 653800
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("F28").Value2 = 0
$ws.Range("G28").Value2 = 0
$ws.Range("H28").Value2 = 0
$ws.Range("I28").Value2 = 0
$ws.Range("J28").Value2 = 0
$ws.Range("K28").Value2 = 0
$ws.Range("D29").Value2 = 0
$ws.Range("E29").Value2 = 0
$ws.Range("F29").Value2 = 0
$ws.Range("G29").Value2 = 0
$ws.Range("H29").Value2 = 0
$ws.Range("I29").Value2 = 0
$ws.Range("J29").Value2 = 0
$ws.Range("K29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("F30").Value2 = 0
$ws.Range("G30").Value2 = 0
$ws.Range("H30").Value2 = 0
$ws.Range("I30").Value2 = 0
$ws.Range("J30").Value2 = 0
$ws.Range("K30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("F31").Value2 = 0
$ws.Range("G31").Value2 = 0
$ws.Range("H31").Value2 = 0
$ws.Range("I31").Value2 = 0
$ws.Range("J31").Value2 = 0
$ws.Range("K31").Value2 = 0
$ws.Range("D32").Value2 = 874000
$ws.Range("E32").Value2 = 2100400
$ws.Range("F32").Value2 = 2134000
$ws.Range("G32").Value2 = -402800
$ws.Range("H32").Value2 = 260300
$ws.Range("I32").Value2 = 462300
$ws.Range("J32").Value2 = 81900
$ws.Range("K32").Value2 = -30500
$ws.Range("D33").Value2 = 2430200
$ws.Range("E33").Value2 = 3883200
$ws.Range("F33").Value2 = 3001300
$ws.Range("G33").Value2 = 3651000
$ws.Range("H33").Value2 = 3280700
$ws.Range("I33").Value2 = 1043500
$ws.Range("J33").Value2 = -6006000
$ws.Range("K33").Value2 = 653800
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("F34").Value2 = 0
$ws.Range("G34").Value2 = 0
$ws.Range("H34").Value2 = 0
$ws.Range("I34").Value2 = 0
$ws.Range("J34").Value2 = 0
$ws.Range("K34").Value2 = 0
$ws.Range("D35").Value2 = 2430200
$ws.Range("E35").Value2 = 3883200
$ws.Range("F35").Value2 = 3001300
$ws.Range("G35").Value2 = 3651000
$ws.Range("H35").Value2 = 3280700
$ws.Range("I35").Value2 = 1043500
$ws.Range("J35").Value2 = -6006000
$ws.Range("K35").Value2 = 653800
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43100
$ws.Range("F38").Value2 = 42735
$ws.Range("G38").Value2 = 42369
$ws.Range("H38").Value2 = 42004
$ws.Range("I38").Value2 = 41639
$ws.Range("J38").Value2 = 41274
$ws.Range("K38").Value2 = 40908
$ws.Range("D41").Value2 = 4127800
$ws.Range("E41").Value2 = 3716000
$ws.Range("F41").Value2 = 8692100
$ws.Range("G41").Value2 = 7738400
$ws.Range("H41").Value2 = 8440700
$ws.Range("I41").Value2 = 8942300
$ws.Range("J41").Value2 = 4517100
$ws.Range("K41").Value2 = 4400600
$ws.Range("D42").Value2 = 13500
$ws.Range("E42").Value2 = 0
$ws.Range("F42").Value2 = 2200
$ws.Range("G42").Value2 = 4500
$ws.Range("H42").Value2 = 3400
$ws.Range("I42").Value2 = 5600
$ws.Range("J42").Value2 = 2351700
$ws.Range("K42").Value2 = 5570900
$ws.Range("D43").Value2 = 16750200
$ws.Range("E43").Value2 = 14556700
$ws.Range("F43").Value2 = 16475300
$ws.Range("G43").Value2 = 13532300
$ws.Range("H43").Value2 = 14776600
$ws.Range("I43").Value2 = 11396100
$ws.Range("J43").Value2 = 9234000
$ws.Range("K43").Value2 = 15544800
$ws.Range("D44").Value2 = 2008400
$ws.Range("E44").Value2 = 2227200
$ws.Range("F44").Value2 = 1827700
$ws.Range("G44").Value2 = 2072300
$ws.Range("H44").Value2 = 1686400
$ws.Range("I44").Value2 = 1191600
$ws.Range("J44").Value2 = 1240900
$ws.Range("K44").Value2 = 2544800
$ws.Range("D45").Value2 = 1638100
$ws.Range("E45").Value2 = 2379700
$ws.Range("F45").Value2 = 2890200
$ws.Range("G45").Value2 = 12762600
$ws.Range("H45").Value2 = 8526000
$ws.Range("I45").Value2 = 3106800
$ws.Range("J45").Value2 = 1773900
$ws.Range("K45").Value2 = 2315900
$ws.Range("D46").Value2 = 24537900
$ws.Range("E46").Value2 = 22879600
$ws.Range("F46").Value2 = 29887600
$ws.Range("G46").Value2 = 36110100
$ws.Range("H46").Value2 = 33433100
$ws.Range("I46").Value2 = 24642300
$ws.Range("J46").Value2 = 16851200
$ws.Range("K46").Value2 = 18622500
$ws.Range("D47").Value2 = 2424600
$ws.Range("E47").Value2 = 7143700
$ws.Range("F47").Value2 = 9661500
$ws.Range("G47").Value2 = 4882900
$ws.Range("H47").Value2 = 3254900
$ws.Range("I47").Value2 = 8447500
$ws.Range("J47").Value2 = 9679400
$ws.Range("K47").Value2 = 10527900
$ws.Range("D48").Value2 = 56807500
$ws.Range("E48").Value2 = 52596600
$ws.Range("F48").Value2 = 52462000
$ws.Range("G48").Value2 = 50082300
$ws.Range("H48").Value2 = 44448800
$ws.Range("I48").Value2 = 41992700
$ws.Range("J48").Value2 = "NA"
$ws.Range("K48").Value2 = 202972000
$ws.Range("D49").Value2 = 72873300
$ws.Range("E49").Value2 = 70533900
$ws.Range("F49").Value2 = 67991500
$ws.Range("G49").Value2 = 63981500
$ws.Range("H49").Value2 = 57855400
$ws.Range("I49").Value2 = 51574500
$ws.Range("J49").Value2 = 77573300
$ws.Range("K49").Value2 = 117761000
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("F50").Value2 = 0
$ws.Range("G50").Value2 = 0
$ws.Range("H50").Value2 = 0
$ws.Range("I50").Value2 = 0
$ws.Range("J50").Value2 = 0
$ws.Range("K50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("F51").Value2 = 0
$ws.Range("G51").Value2 = 0
$ws.Range("H51").Value2 = 0
$ws.Range("I51").Value2 = 0
$ws.Range("J51").Value2 = 0
$ws.Range("K51").Value2 = 0
$ws.Range("D52").Value2 = 6466000
$ws.Range("E52").Value2 = 5421500
$ws.Range("F52").Value2 = 6596200
$ws.Range("G52").Value2 = 6420000
$ws.Range("H52").Value2 = 6148500
$ws.Range("I52").Value2 = 5903900
$ws.Range("J52").Value2 = 7243600
$ws.Range("K52").Value2 = 9460900
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("F53").Value2 = 0
$ws.Range("G53").Value2 = 0
$ws.Range("H53").Value2 = 0
$ws.Range("I53").Value2 = 0
$ws.Range("J53").Value2 = 0
$ws.Range("K53").Value2 = 0
$ws.Range("D54").Value2 = 163109000
$ws.Range("E54").Value2 = 158575000
$ws.Range("F54").Value2 = 166599000
$ws.Range("G54").Value2 = 161477000
$ws.Range("H54").Value2 = 145141000
$ws.Range("I54").Value2 = 132561000
$ws.Range("J54").Value2 = 121110000
$ws.Range("K54").Value2 = 143787000
$ws.Range("D57").Value2 = 12004200
$ws.Range("E57").Value2 = 12267800
$ws.Range("F57").Value2 = 11655200
$ws.Range("G57").Value2 = 12383400
$ws.Range("H57").Value2 = 10805900
$ws.Range("I57").Value2 = 8113100
$ws.Range("J57").Value2 = 7197600
$ws.Range("K57").Value2 = 15077600
$ws.Range("D58").Value2 = 8715600
$ws.Range("E58").Value2 = 6957500
$ws.Range("F58").Value2 = 12085000
$ws.Range("G58").Value2 = 13754500
$ws.Range("H58").Value2 = 9009600
$ws.Range("I58").Value2 = 6440200
$ws.Range("J58").Value2 = 7803400
$ws.Range("K58").Value2 = 23990300
$ws.Range("D59").Value2 = 11979500
$ws.Range("E59").Value2 = 11479100
$ws.Range("F59").Value2 = 13426900
$ws.Range("G59").Value2 = 11502600
$ws.Range("H59").Value2 = 11822400
$ws.Range("I59").Value2 = 10687000
$ws.Range("J59").Value2 = 18499400
$ws.Range("K59").Value2 = 17246800
$ws.Range("D60").Value2 = 32699300
$ws.Range("E60").Value2 = 30704400
$ws.Range("F60").Value2 = 37167000
$ws.Range("G60").Value2 = 37640500
$ws.Range("H60").Value2 = 31637900
$ws.Range("I60").Value2 = 25240300
$ws.Range("J60").Value2 = 25800200
$ws.Range("K60").Value2 = 28423800
$ws.Range("D61").Value2 = 57190100
$ws.Range("E61").Value2 = 54214600
$ws.Range("F61").Value2 = 55843700
$ws.Range("G61").Value2 = 52947800
$ws.Range("H61").Value2 = 49345100
$ws.Range("I61").Value2 = 47414200
$ws.Range("J61").Value2 = 38769200
$ws.Range("K61").Value2 = 44721000
$ws.Range("D62").Value2 = 24484100
$ws.Range("E62").Value2 = 26005500
$ws.Range("F62").Value2 = 30004300
$ws.Range("G62").Value2 = 28084500
$ws.Range("H62").Value2 = 25935900
$ws.Range("I62").Value2 = 23932000
$ws.Range("J62").Value2 = 28755500
$ws.Range("K62").Value2 = 30198600
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("F63").Value2 = 0
$ws.Range("G63").Value2 = 0
$ws.Range("H63").Value2 = 0
$ws.Range("I63").Value2 = 0
$ws.Range("J63").Value2 = 0
$ws.Range("K63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("F64").Value2 = 0
$ws.Range("G64").Value2 = 0
$ws.Range("H64").Value2 = 0
$ws.Range("I64").Value2 = 0
$ws.Range("J64").Value2 = 0
$ws.Range("K64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("F65").Value2 = 0
$ws.Range("G65").Value2 = 0
$ws.Range("H65").Value2 = 0
$ws.Range("I65").Value2 = 0
$ws.Range("J65").Value2 = 0
$ws.Range("K65").Value2 = 0
$ws.Range("D66").Value2 = 128432000
$ws.Range("E66").Value2 = 124093000
$ws.Range("F66").Value2 = 133719000
$ws.Range("G66").Value2 = 128490000
$ws.Range("H66").Value2 = 116601000
$ws.Range("I66").Value2 = 105769000
$ws.Range("J66").Value2 = 92018900
$ws.Range("K66").Value2 = 102232000
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("F68").Value2 = 0
$ws.Range("G68").Value2 = 0
$ws.Range("H68").Value2 = 0
$ws.Range("I68").Value2 = 0
$ws.Range("J68").Value2 = 0
$ws.Range("K68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("F69").Value2 = 0
$ws.Range("G69").Value2 = 0
$ws.Range("H69").Value2 = 0
$ws.Range("I69").Value2 = 0
$ws.Range("J69").Value2 = 0
$ws.Range("K69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("F70").Value2 = 0
$ws.Range("G70").Value2 = 0
$ws.Range("H70").Value2 = 0
$ws.Range("I70").Value2 = 0
$ws.Range("J70").Value2 = 0
$ws.Range("K70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("F71").Value2 = 0
$ws.Range("G71").Value2 = 0
$ws.Range("H71").Value2 = 0
$ws.Range("I71").Value2 = 0
$ws.Range("J71").Value2 = 0
$ws.Range("K71").Value2 = 0
$ws.Range("D72").Value2 = -39523200
$ws.Range("E72").Value2 = -39593900
$ws.Range("F72").Value2 = -40450000
$ws.Range("G72").Value2 = -40071900
$ws.Range("H72").Value2 = -41355400
$ws.Range("I72").Value2 = -40960500
$ws.Range("J72").Value2 = -77334300
$ws.Range("K72").Value2 = -58425200
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("F73").Value2 = 0
$ws.Range("G73").Value2 = 0
$ws.Range("H73").Value2 = 0
$ws.Range("I73").Value2 = 0
$ws.Range("J73").Value2 = 0
$ws.Range("K73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("F74").Value2 = 0
$ws.Range("G74").Value2 = 0
$ws.Range("H74").Value2 = 0
$ws.Range("I74").Value2 = 0
$ws.Range("J74").Value2 = 0
$ws.Range("K74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("F75").Value2 = 0
$ws.Range("G75").Value2 = 0
$ws.Range("H75").Value2 = 0
$ws.Range("I75").Value2 = 0
$ws.Range("J75").Value2 = 0
$ws.Range("K75").Value2 = 0
$ws.Range("D76").Value2 = 34677300
$ws.Range("E76").Value2 = 34482100
$ws.Range("F76").Value2 = 32879900
$ws.Range("G76").Value2 = 32986500
$ws.Range("H76").Value2 = 28540100
$ws.Range("I76").Value2 = 26792000
$ws.Range("J76").Value2 = 29091000
$ws.Range("K76").Value2 = 41555200
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("F77").Value2 = 0
$ws.Range("G77").Value2 = 0
$ws.Range("H77").Value2 = 0
$ws.Range("I77").Value2 = 0
$ws.Range("J77").Value2 = 0
$ws.Range("K77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43100
$ws.Range("F80").Value2 = 42735
$ws.Range("G80").Value2 = 42369
$ws.Range("H80").Value2 = 42004
$ws.Range("I80").Value2 = 41639
$ws.Range("J80").Value2 = 41274
$ws.Range("K80").Value2 = 40908
$ws.Range("D81").Value2 = 2430200
$ws.Range("E81").Value2 = 3883200
$ws.Range("F81").Value2 = 3001300
$ws.Range("G81").Value2 = 3651000
$ws.Range("H81").Value2 = 3280700
$ws.Range("I81").Value2 = 1043500
$ws.Range("J81").Value2 = -6006000
$ws.Range("K81").Value2 = 653800
$ws.Range("D83").Value2 = 15523900
$ws.Range("E83").Value2 = 16365300
$ws.Range("F83").Value2 = 15012200
$ws.Range("G83").Value2 = 12745800
$ws.Range("H83").Value2 = 11863900
$ws.Range("I83").Value2 = 12234200
$ws.Range("J83").Value2 = 24635500
$ws.Range("K83").Value2 = 16945100
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("F84").Value2 = 0
$ws.Range("G84").Value2 = 0
$ws.Range("H84").Value2 = 0
$ws.Range("I84").Value2 = 0
$ws.Range("J84").Value2 = 0
$ws.Range("K84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("F85").Value2 = 0
$ws.Range("G85").Value2 = 0
$ws.Range("H85").Value2 = 0
$ws.Range("I85").Value2 = 0
$ws.Range("J85").Value2 = 0
$ws.Range("K85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("F86").Value2 = 0
$ws.Range("G86").Value2 = 0
$ws.Range("H86").Value2 = 0
$ws.Range("I86").Value2 = 0
$ws.Range("J86").Value2 = 0
$ws.Range("K86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("F87").Value2 = 0
$ws.Range("G87").Value2 = 0
$ws.Range("H87").Value2 = 0
$ws.Range("I87").Value2 = 0
$ws.Range("J87").Value2 = 0
$ws.Range("K87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("F88").Value2 = 0
$ws.Range("G88").Value2 = 0
$ws.Range("H88").Value2 = 0
$ws.Range("I88").Value2 = 0
$ws.Range("J88").Value2 = 0
$ws.Range("K88").Value2 = 0
$ws.Range("D89").Value2 = 20128500
$ws.Range("E89").Value2 = 19297100
$ws.Range("F89").Value2 = 17424500
$ws.Range("G89").Value2 = 16827600
$ws.Range("H89").Value2 = 15062700
$ws.Range("I89").Value2 = 14569000
$ws.Range("J89").Value2 = 15233300
$ws.Range("K89").Value2 = 19032200
$ws.Range("D91").Value2 = -10253900
$ws.Range("E91").Value2 = -10265100
$ws.Range("F91").Value2 = -9017400
$ws.Range("G91").Value2 = -9163300
$ws.Range("H91").Value2 = -8062600
$ws.Range("I91").Value2 = -7371500
$ws.Range("J91").Value2 = -6306700
$ws.Range("K91").Value2 = -9867000
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("F92").Value2 = 0
$ws.Range("G92").Value2 = 0
$ws.Range("H92").Value2 = 0
$ws.Range("I92").Value2 = 0
$ws.Range("J92").Value2 = 0
$ws.Range("K92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("F93").Value2 = 0
$ws.Range("G93").Value2 = 0
$ws.Range("H93").Value2 = 0
$ws.Range("I93").Value2 = 0
$ws.Range("J93").Value2 = 0
$ws.Range("K93").Value2 = 0
$ws.Range("D94").Value2 = -16041100
$ws.Range("E94").Value2 = -18865100
$ws.Range("F94").Value2 = -15268000
$ws.Range("G94").Value2 = -16846700
$ws.Range("H94").Value2 = -12073700
$ws.Range("I94").Value2 = -11103200
$ws.Range("J94").Value2 = -7484800
$ws.Range("K94").Value2 = -10887100
$ws.Range("D96").Value2 = -3651000
$ws.Range("E96").Value2 = -1749200
$ws.Range("F96").Value2 = -1790700
$ws.Range("G96").Value2 = -1409200
$ws.Range("H96").Value2 = -1447400
$ws.Range("I96").Value2 = -2516600
$ws.Range("J96").Value2 = -3814800
$ws.Range("K96").Value2 = -4133000
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("F97").Value2 = 0
$ws.Range("G97").Value2 = 0
$ws.Range("H97").Value2 = 0
$ws.Range("I97").Value2 = 0
$ws.Range("J97").Value2 = 0
$ws.Range("K97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("F98").Value2 = 0
$ws.Range("G98").Value2 = 0
$ws.Range("H98").Value2 = 0
$ws.Range("I98").Value2 = 0
$ws.Range("J98").Value2 = 0
$ws.Range("K98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("F99").Value2 = 0
$ws.Range("G99").Value2 = 0
$ws.Range("H99").Value2 = 0
$ws.Range("I99").Value2 = 0
$ws.Range("J99").Value2 = 0
$ws.Range("K99").Value2 = 0
$ws.Range("D100").Value2 = -3656600
$ws.Range("E100").Value2 = -5154400
$ws.Range("F100").Value2 = -1483300
$ws.Range("G100").Value2 = -982900
$ws.Range("H100").Value2 = -3852900
$ws.Range("I100").Value2 = 1146700
$ws.Range("J100").Value2 = -7406300
$ws.Range("K100").Value2 = -6993600
$ws.Range("D101").Value2 = -19100
$ws.Range("E101").Value2 = -253600
$ws.Range("F101").Value2 = 280500
$ws.Range("G101").Value2 = 299600
$ws.Range("H101").Value2 = 362400
$ws.Range("I101").Value2 = -187400
$ws.Range("J101").Value2 = -31400
$ws.Range("K101").Value2 = -47000
$ws.Range("D102").Value2 = 411800
$ws.Range("E102").Value2 = -4976000
$ws.Range("F102").Value2 = 953700
$ws.Range("G102").Value2 = -702400
$ws.Range("H102").Value2 = -501500
$ws.Range("I102").Value2 = 4425100
$ws.Range("J102").Value2 = 310800
$ws.Range("K102").Value2 = 1104600
